$d = $word.ActiveDocument

# 1. Name casing: DHEERAJ CHAND -> Dheeraj Chand
$d.Content.Find.Execute("DHEERAJ CHAND", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dheeraj Chand", 2) | Out-Null

# 2. Professional title placeholder
$d.Content.Find.Execute("Director of Research and Analysis", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Professional Title", 2) | Out-Null

# 3. Contact line: phone/email format update
$d.Content.Find.Execute("(202) 550-7110 | Dheeraj.Chand@gmail.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "202.550.7110 | dheeraj.chand@gmail.com", 2) | Out-Null

# 4. Years of experience: 20+ -> 21
$d.Content.Find.Execute("20+ years of experience", $true, $false, $false, $false, $false,
                         $true, 1, $false, "21 years of experience", 2) | Out-Null

# 5. Company name placeholder for the PARTNER role
$d.Content.Find.Execute("Siege Analytics, Austin, TX | 2005", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Your Company Name, Your City, ST | 2005", 2) | Out-Null

# 6. Remove the DATA PRODUCTS MANAGER, ANALYTICS SUPERVISOR, SOFTWARE ENGINEER,
#    and RESEARCH DIRECTOR job blocks entirely (everything between the PARTNER
#    bullet list and the KEY ACHIEVEMENTS AND IMPACT heading).
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "DATA PRODUCTS MANAGER") {
        $startPara = $i
    }
    if ($t -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $endPara = $i - 1
        break
    }
}
if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($d.Paragraphs($startPara).Range.Start, $d.Paragraphs($endPara).Range.End)
    $r.Delete()
}

# 7. Remove the "Systems and Infrastructure Development" and
#    "Community and Stakeholder Engagement" sections from KEY ACHIEVEMENTS,
#    keeping only "Research Leadership and Community Impact".
$startPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Systems and Infrastructure Development") {
        $startPara2 = $i
        break
    }
}
$endPara2 = $d.Paragraphs.Count
if ($startPara2 -ne $null) {
    $r2 = $d.Range($d.Paragraphs($startPara2).Range.Start, $d.Paragraphs($endPara2).Range.End)
    $r2.Delete()
}
